# Change highlight color to cyan for the paragraphs describing
# "Landlord, Agent and Prospects Receive and Send Messages":
#   - "View Prospects Notifications" (was lightGray -> cyan)
#   - "Reply Prospect" (no highlight -> cyan)
#   - "Send Property Inquiry to Landlord or Agent" (no highlight -> cyan)
#
# wdTurquoise (3) is the WdColorIndex value that Word serialises as
# w:highlight w:val="cyan" in the OOXML.

$d = $word.ActiveDocument
$wdTurquoise = 3

$targets = @(
    "View Prospects Notifications",
    "Reply Prospect",
    "Send Property Inquiry to Landlord or Agent"
)

foreach ($target in $targets) {
    $rng = $d.Content
    $found = $rng.Find.Execute($target, $false, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $para = $rng.Paragraphs(1)
        # Use the paragraph's full range (which includes the trailing
        # paragraph mark) so the highlight is applied both to the run(s)
        # and to the paragraph mark's run properties (w:pPr/w:rPr).
        $para.Range.Font.HighlightColorIndex = $wdTurquoise
    }
}
